# Rename the BTEC / Pearson logo pictures embedded in the headers and
# footers of every section: the two copies of the BTec logo become
# "image1.jpg" (they were "image2.jpg"), and the two copies of the
# Pearson Edexcel logo become "image2.png" (they were "image1.png").
#
# InlineShape.Name does not read back the name already stored in the
# document (it only reflects names set earlier in this same script),
# so shapes are identified by their AlternativeText (the picture's
# description), which does reflect the document as loaded.

$d = $word.ActiveDocument

For ($secIdx = 1; $secIdx -le $d.Sections.Count; $secIdx++) {
    $sec = $d.Sections.Item($secIdx)

    For ($hIdx = 1; $hIdx -le 3; $hIdx++) {
        $hdr = $sec.Headers.Item($hIdx)
        If ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            For ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                If ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }

    For ($fIdx = 1; $fIdx -le 3; $fIdx++) {
        $ftr = $sec.Footers.Item($fIdx)
        If ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            For ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                If ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}
